# Update "想去人数" (want-to-go count, column F) figures to the values
# captured at the later gh-pages data-refresh commit (456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (exhibitions)
$wsExhibit.Cells.Item(7, 6).Value  = 6290
$wsExhibit.Cells.Item(12, 6).Value = 9492
$wsExhibit.Cells.Item(14, 6).Value = 2565
$wsExhibit.Cells.Item(16, 6).Value = 2357
$wsExhibit.Cells.Item(17, 6).Value = 2557
$wsExhibit.Cells.Item(28, 6).Value = 113
$wsExhibit.Cells.Item(29, 6).Value = 595
$wsExhibit.Cells.Item(30, 6).Value = 1254
$wsExhibit.Cells.Item(31, 6).Value = 1231
$wsExhibit.Cells.Item(35, 6).Value = 1609
$wsExhibit.Cells.Item(36, 6).Value = 2663
$wsExhibit.Cells.Item(41, 6).Value = 24

# Sheet "本地生活" (local life)
$wsLocal.Cells.Item(3, 6).Value = 927

# Sheet "全部类型" (all types, combined roll-up)
$wsAll.Cells.Item(4, 6).Value  = 927
$wsAll.Cells.Item(10, 6).Value = 6290
$wsAll.Cells.Item(14, 6).Value = 9492
$wsAll.Cells.Item(17, 6).Value = 2565
$wsAll.Cells.Item(19, 6).Value = 2357
$wsAll.Cells.Item(20, 6).Value = 2557
$wsAll.Cells.Item(30, 6).Value = 113
$wsAll.Cells.Item(31, 6).Value = 595
$wsAll.Cells.Item(32, 6).Value = 1254
$wsAll.Cells.Item(33, 6).Value = 1231
$wsAll.Cells.Item(37, 6).Value = 1609
$wsAll.Cells.Item(39, 6).Value = 2663
